$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '42.052.23'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -0.47%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.217.79'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -1.33%  '
$ws.Range('E4').Value = '  +0.12%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '242.21'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -1.80%  '
$ws.Range('E6').Value = '  +0.87%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '73.25'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -1.88%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.607'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -2.08%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '42.33'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -0.91%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0956'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +0.77%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '7.07'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -0.88%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.104'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +0.78%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '2.550.70'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -1.32%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '14.29'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -1.41%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.837'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -1.58%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '2.211.10'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -2.54%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '41.895.33'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -0.57%  '
$ws.Range('E19').Value = '  +5.87%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '6.20'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +1.04%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '72.96'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +0.73%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '10.61'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +18.55%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '230.54'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -0.51%  '
$ws.Range('E24').Value = '  -6.28%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '11.87'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +3.78%  '
$ws.Range('E26').Value = '  +0.08%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '3.66'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +1.30%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.28'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -1.57%  '
$ws.Range('E29').Value = '  -0.08%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '168.13'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -0.34%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '20.51'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -0.95%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '5.65'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +7.63%  '
$ws.Range('E33').Value = '  -2.73%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '29.85'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -3.20%  '
$ws.Range('E35').Value = '  -0.15%  '
$ws.Range('E36').Value = '  -10.03%  '
$ws.Range('E37').Value = '  -4.14%  '
$ws.Range('E38').Value = '  -4.11%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '13.81'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +0.15%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '65.91'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +5.32%  '
$ws.Range('E41').Value = '  -2.18%  '
$ws.Range('E42').Value = '  -2.20%  '
$ws.Range('E43').Value = '  -2.70%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '8.81'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +1.44%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '105.54'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -1.19%  '
$ws.Range('E46').Value = '  -1.90%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.48'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +7.69%  '
$ws.Range('E48').Value = '  -0.36%  '
$ws.Range('E49').Value = '  -0.61%  '
$ws.Range('E50').Value = '  -0.17%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '2.423.88'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -1.45%  '
